$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.317.78"
$ws.Cells.Item(2, 5).Value = "  +0.42%  "

$ws.Cells.Item(3, 4).Value = "3.169.83"
$ws.Cells.Item(3, 5).Value = "  -1.22%  "

$ws.Cells.Item(4, 5).Value = "  -0.03%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "602.92"
$ws.Cells.Item(5, 5).Value = "  -0.08%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "154.01"
$ws.Cells.Item(6, 5).Value = "  +0.00%  "

$ws.Cells.Item(7, 5).Value = "  +0.01%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.550"
$ws.Cells.Item(8, 5).Value = "  +2.81%  "

$ws.Cells.Item(9, 4).Value = "3.169.45"

$ws.Cells.Item(10, 5).Value = "  -2.12%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.57"
$ws.Cells.Item(11, 5).Value = "  -9.71%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.510"
$ws.Cells.Item(12, 5).Value = "  +0.21%  "

$ws.Cells.Item(13, 5).Value = "  -3.27%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "38.50"
$ws.Cells.Item(14, 5).Value = "  -0.80%  "

$ws.Cells.Item(15, 4).Value = "3.692.03"
$ws.Cells.Item(15, 5).Value = "  -1.33%  "

$ws.Cells.Item(16, 4).Value = "66.358.24"
$ws.Cells.Item(16, 5).Value = "  +0.22%  "

$ws.Cells.Item(17, 5).Value = "  -0.85%  "

$ws.Cells.Item(18, 4).Value = "3.174.61"
$ws.Cells.Item(18, 5).Value = "  -1.44%  "

$ws.Cells.Item(19, 5).Value = "  +0.35%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "509.74"
$ws.Cells.Item(20, 5).Value = "  -0.19%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "15.39"
$ws.Cells.Item(21, 5).Value = "  -1.54%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.729"
$ws.Cells.Item(22, 5).Value = "  -0.96%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.09"
$ws.Cells.Item(23, 5).Value = "  +1.44%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "14.74"
$ws.Cells.Item(24, 5).Value = "  -2.98%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "84.65"
$ws.Cells.Item(25, 5).Value = "  -0.82%  "

$ws.Cells.Item(26, 5).Value = "  -0.07%  "

$ws.Cells.Item(27, 5).Value = "  -1.08%  "

$ws.Cells.Item(28, 5).Value = "  -1.36%  "

$ws.Cells.Item(29, 5).Value = "  +6.67%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "3.06"
$ws.Cells.Item(30, 5).Value = "  +6.32%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.91"
$ws.Cells.Item(31, 5).Value = "  +0.14%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "27.93"
$ws.Cells.Item(32, 5).Value = "  -0.97%  "

$ws.Cells.Item(33, 5).Value = "  -0.16%  "

$ws.Cells.Item(34, 5).Value = "  -1.46%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "6.50"
$ws.Cells.Item(35, 5).Value = "  -1.91%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "510.62"
$ws.Cells.Item(36, 5).Value = "  +5.52%  "

$ws.Cells.Item(37, 5).Value = "  -1.05%  "

$ws.Cells.Item(38, 5).Value = "  -3.05%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0418"

$ws.Cells.Item(40, 5).Value = "  +6.76%  "

$ws.Cells.Item(41, 5).Value = "  -0.22%  "

$ws.Cells.Item(42, 5).Value = "  +4.48%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.296"
$ws.Cells.Item(43, 5).Value = "  -0.14%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.78"
$ws.Cells.Item(44, 5).Value = "  -7.29%  "

$ws.Cells.Item(45, 5).Value = "  -2.97%  "

$ws.Cells.Item(46, 4).Value = "2.833.88"
$ws.Cells.Item(46, 5).Value = "  -4.32%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "27.97"
$ws.Cells.Item(47, 5).Value = "  -3.37%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.36"
$ws.Cells.Item(49, 5).Value = "  +1.01%  "

$ws.Cells.Item(50, 5).Value = "  +0.52%  "

$ws.Cells.Item(51, 2).Value = "Arweave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "34.74"
$ws.Cells.Item(51, 5).Value = "  +2.13%  "
